$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) date column C was updated for every data row
# (rows 2-372) from serial date 45190 (2023-09-21) to 45192 (2023-09-23).
$ws.Range("C2:C372").Value = 45192
